# Apply "results with fixed workflow" update:
# - Remove the last 4 data rows (rows 17-20) from both sheets, shrinking the
#   table from 19 data points down to 15.
# - Shift the remaining "Cutoff" (column B) values by +4 (e.g. 1->5, 2->6, ...)
# - Update the "Reaction_number" (column C) values to the refreshed results.

$wb = $excel.ActiveWorkbook

$sheetData = @{
    "NBR" = @(860, 857, 856, 856, 849, 847, 844, 840, 826, 824, 814, 805, 800, 794, 794)
    "BAR" = @(715, 708, 706, 699, 695, 690, 684, 684, 682, 682, 678, 679, 682, 680, 673)
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $sheetData.ContainsKey($name)) {
        continue
    }

    # Drop the trailing 4 rows (old rows 17-20) so the table ends at row 16.
    $ws.Range("A17:C20").EntireRow.Delete() | Out-Null

    $values = $sheetData[$name]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $i + 5
        $ws.Cells.Item($row, 3).Value = $values[$i]
    }
}
